$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add I1 = "I0" and J1 = "IF", matching the style of the other header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-39: I column is constant 1
$ws.Range("I2:I39").Value = 1

# Data rows 2-39: J column duplicates the H column values
$ws.Range("H2:H39").Copy()
$ws.Range("J2:J39").PasteSpecial(-4163)
